$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '62.661.45'
$ws.Range('E2').Value = '  +3.36%  '
$ws.Range('D3').Value = '2.445.86'
$ws.Range('E3').Value = '  +2.01%  '
$ws.Range('E4').Value = '  -0.16%  '
Set-TextValue $ws 'D5' '577.54'
$ws.Range('E5').Value = '  +2.47%  '
Set-TextValue $ws 'D6' '145.68'
$ws.Range('E6').Value = '  +3.35%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').Value = '2.444.78'
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('E10').Value = '  +2.57%  '
$ws.Range('E11').Value = '  +0.86%  '
Set-TextValue $ws 'D12' '5.22'
$ws.Range('E12').Value = '  +1.03%  '
Set-TextValue $ws 'D13' '0.353'
$ws.Range('E13').Value = '  +3.32%  '
Set-TextValue $ws 'D14' '28.47'
$ws.Range('E14').Value = '  +9.12%  '
Set-TextValue $ws 'D15' '0.0000178'
$ws.Range('E15').Value = '  +5.87%  '
$ws.Range('D16').Value = '2.883.80'
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('D17').Value = '62.530.35'
$ws.Range('E17').Value = '  +3.15%  '
$ws.Range('D18').Value = '2.433.73'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws 'D19' '7.81'
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D20' '10.91'
$ws.Range('E20').Value = '  +2.72%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws 'D21' '325.98'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws 'D22' '4.13'
$ws.Range('E22').Value = '  +1.07%  '
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws 'D23' '2.02'
$ws.Range('E23').Value = '  +10.82%  '
$ws.Range('B24').Value = 'BabyDogeCoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D24').Value = '0.0₆0704'
$ws.Range('E24').Value = '  +148.02%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  +0.42%  '
Set-TextValue $ws 'D27' '652.87'
$ws.Range('E27').Value = '  +15.43%  '
$ws.Range('E28').Value = '  +14.04%  '
Set-TextValue $ws 'D29' '8.55'
$ws.Range('E29').Value = '  +6.10%  '
$ws.Range('E30').Value = '  +4.69%  '
$ws.Range('D31').Value = '2.558.67'
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('E32').Value = '  +1.39%  '
Set-TextValue $ws 'D33' '1.42'
$ws.Range('E33').Value = '  +6.26%  '
Set-TextValue $ws 'D34' '1.87'
$ws.Range('E34').Value = '  +2.98%  '
$ws.Range('E35').Value = '  +5.80%  '
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +3.37%  '
Set-TextValue $ws 'D39' '5.48'
$ws.Range('E39').Value = '  +6.57%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws 'D40' '0.373'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D41' '152.76'
$ws.Range('E41').Value = '  +0.30%  '
Set-TextValue $ws 'D42' '18.56'
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('E43').Value = '  +9.00%  '
$ws.Range('E44').Value = '  +4.94%  '
Set-TextValue $ws 'D45' '42.54'
$ws.Range('E45').Value = '  +2.08%  '
Set-TextValue $ws 'D46' '0.999'
$ws.Range('E46').Value = '  +0.03%  '
Set-TextValue $ws 'D47' '15.01'
Set-TextValue $ws 'D48' '144.19'
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('E49').Value = '  +1.08%  '
Set-TextValue $ws 'D50' '20.53'
$ws.Range('E50').Value = '  +6.34%  '
Set-TextValue $ws 'D51' '0.604'
$ws.Range('E51').Value = '  +2.35%  '
